$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix "Marking" row (row 11): Right = 4, Wrong = -2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Fix "Total" row (row 12): Right = 60, Wrong = -2, Max display = 58 / 112
$ws.Range("B12").Value = 60
$ws.Range("C12").Value = -2
$ws.Range("E12").Value = "58 / 112"
